# Update the "Förändrad" (Changed) date column (C) for rows 2-14
# from serial date 45184 (2023-09-15) to serial date 45186 (2023-09-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}
